# "bikin biar ga ada duplikat masuk di tamu" - append 5 new, unique guest
# entries to the Tamu (guest) list so no duplicates are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Append the new guest rows (27-31) under the existing data (2-26)
# ---------------------------------------------------------------------
$newGuests = @(
    @("Lukman Hakim",      "lukman.hakim@travelindo.com",          "PT Travel Indo Raya"),
    @("Dewi Sartika",      "dewi.sartika@poltek-x.ac.id",          "Politeknik Negeri X"),
    @("Rio Ferdinand",     "rio.ferdinand@kesehatandaerah.go.id",  "Dinas Kesehatan Provinsi"),
    @("Cindy Aurelia",     "cindy.aurelia@fashionista.net",        "Aurelia Fashion Group"),
    @("Haikal Zulkarnain", "haikal.zul@softwareid.com",            "ID Software Development")
)

$startRow = 27
$lastRow = $startRow + $newGuests.Count - 1

# Carry the header row's look (bold font, medium border) onto the new rows,
# matching the style reused for rows 27-31 in the saved workbook.
$ws.Range("A1:C1").Copy() | Out-Null
$ws.Range("A" + $startRow + ":C" + $lastRow).PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $newGuests.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newGuests[$i][0]
    $ws.Cells.Item($r, 2).Value = $newGuests[$i][1]
    $ws.Cells.Item($r, 3).Value = $newGuests[$i][2]
    $ws.Rows.Item($r).RowHeight = 15
}

# ---------------------------------------------------------------------
# 2. Re-centre the header / data text (was left-aligned, now centred)
# ---------------------------------------------------------------------
$ws.Range("A1:C" + $lastRow).HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Resize columns to fit the new, wider content
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.666666
$ws.Columns.Item(2).ColumnWidth = 39.666666
$ws.Columns.Item(3).ColumnWidth = 41.166666

# ---------------------------------------------------------------------
# 4. Narrow the saved selection to column A only (full column)
# ---------------------------------------------------------------------
$ws.Range("A1:A1048576").Select() | Out-Null
